$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 6565
$ws.Range("E2").Value = 70
$ws.Range("F2").Value = 70
$ws.Range("G2").Value = 11
$ws.Range("H2").Value = 26
$ws.Range("I2").Value = $null
$ws.Range("J2").Value = 26
$ws.Range("K2").Value = 4698
$ws.Range("L2").Value = 2893
$ws.Range("M2").Value = 1805
$ws.Range("N2").Value = $null
$ws.Range("O2").Value = 1805
$ws.Range("P2").Value = 1637
$ws.Range("Q2").Value = 177
$ws.Range("R2").Value = -388
$ws.Range("S2").Value = 222
$ws.Range("T2").Value = 367
$ws.Range("U2").Value = -190
$ws.Range("V2").Value = 1602
$ws.Range("W2").Value = 1.07
$ws.Range("X2").Value = 0.4
$ws.Range("Y2").Value = 1.47
$ws.Range("Z2").Value = 0.57
$ws.Range("AA2").Value = 160.22
$ws.Range("AB2").Value = 10.97
$ws.Range("AC2").Value = 81
$ws.Range("AD2").Value = 57.14
$ws.Range("AE2").Value = 5516
$ws.Range("AF2").Value = 0.83
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 24367748

# Row 3
$ws.Range("D3").Value = 6786
$ws.Range("E3").Value = 40
$ws.Range("F3").Value = 40
$ws.Range("G3").Value = -28
$ws.Range("H3").Value = -14
$ws.Range("I3").Value = -15
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 5191
$ws.Range("L3").Value = 3150
$ws.Range("M3").Value = 2041
$ws.Range("N3").Value = 2010
$ws.Range("O3").Value = 31
$ws.Range("P3").Value = 1880
$ws.Range("Q3").Value = 163
$ws.Range("R3").Value = -587
$ws.Range("S3").Value = 434
$ws.Range("T3").Value = 573
$ws.Range("U3").Value = -409
$ws.Range("V3").Value = 1822
$ws.Range("W3").Value = 0.58
$ws.Range("X3").Value = -0.21
$ws.Range("Y3").Value = -0.79
$ws.Range("Z3").Value = -0.29
$ws.Range("AA3").Value = 154.34
$ws.Range("AB3").Value = 6.94
$ws.Range("AC3").Value = -42
$ws.Range("AD3").Value = -128.02
$ws.Range("AE3").Value = 5348
$ws.Range("AF3").Value = 1.01
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 32228965

# Row 4
$ws.Range("D4").Value = 7060
$ws.Range("E4").Value = 183
$ws.Range("F4").Value = 183
$ws.Range("G4").Value = 104
$ws.Range("H4").Value = 76
$ws.Range("I4").Value = 75
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 5225
$ws.Range("L4").Value = 3118
$ws.Range("M4").Value = 2106
$ws.Range("N4").Value = 2086
$ws.Range("O4").Value = 20
$ws.Range("P4").Value = 1880
$ws.Range("Q4").Value = 501
$ws.Range("R4").Value = -251
$ws.Range("S4").Value = -217
$ws.Range("T4").Value = 281
$ws.Range("U4").Value = 219
$ws.Range("V4").Value = 1646
$ws.Range("W4").Value = 2.6
$ws.Range("X4").Value = 1.08
$ws.Range("Y4").Value = 3.65
$ws.Range("Z4").Value = 1.47
$ws.Range("AA4").Value = 148.05
$ws.Range("AB4").Value = 10.94
$ws.Range("AC4").Value = 199
$ws.Range("AD4").Value = 28.81
$ws.Range("AE4").Value = 5549
$ws.Range("AF4").Value = 1.03
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 37233066

# Row 5
$ws.Range("D5").Value = 6599
$ws.Range("E5").Value = -253
$ws.Range("F5").Value = -253
$ws.Range("G5").Value = -235
$ws.Range("H5").Value = -222
$ws.Range("I5").Value = -225
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 5731
$ws.Range("L5").Value = 3850
$ws.Range("M5").Value = 1881
$ws.Range("N5").Value = 1861
$ws.Range("O5").Value = 20
$ws.Range("P5").Value = 1880
$ws.Range("Q5").Value = -29
$ws.Range("R5").Value = -232
$ws.Range("S5").Value = 258
$ws.Range("T5").Value = 370
$ws.Range("U5").Value = -399
$ws.Range("V5").Value = 2028
$ws.Range("W5").Value = -3.83
$ws.Range("X5").Value = -3.36
$ws.Range("Y5").Value = -11.4
$ws.Range("Z5").Value = -4.05
$ws.Range("AA5").Value = 204.65
$ws.Range("AB5").Value = -0.92
$ws.Range("AC5").Value = -598
$ws.Range("AD5").Value = -6.95
$ws.Range("AE5").Value = 4949
$ws.Range("AF5").Value = 0.84
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 37240693

# Row 6
$ws.Range("D6").Value = 6263
$ws.Range("E6").Value = -292
$ws.Range("F6").Value = -292
$ws.Range("G6").Value = -317
$ws.Range("H6").Value = -336
$ws.Range("I6").Value = -338
$ws.Range("K6").Value = 5994
$ws.Range("L6").Value = 4247
$ws.Range("M6").Value = 1747
$ws.Range("N6").Value = 1727
$ws.Range("P6").Value = 1880
$ws.Range("Q6").Value = -225
$ws.Range("R6").Value = -142
$ws.Range("S6").Value = 370
$ws.Range("T6").Value = 423
$ws.Range("U6").Value = -647
$ws.Range("V6").Value = 2582
$ws.Range("W6").Value = -4.67
$ws.Range("X6").Value = -5.37
$ws.Range("Y6").Value = -18.84
$ws.Range("Z6").Value = -5.74
$ws.Range("AA6").Value = 243.17
$ws.Range("AB6").Value = -19.01
$ws.Range("AC6").Value = -899
$ws.Range("AD6").Value = -3.39
$ws.Range("AE6").Value = 4594
$ws.Range("AF6").Value = 0.66
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 37240693

# Row 7
$ws.Range("D7:AI7").Value = $null

# Row 8
$ws.Range("D8:AI8").Value = $null

# Row 9
$ws.Range("D9:AI9").Value = $null
